$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1844660194174757
$ws.Range("C2").Value = 0.5598705501618123
$ws.Range("J2").Value = 0.003236245954692557
$ws.Range("P2").Value = 0.116504854368932
$ws.Range("S2").Value = 0.1359223300970874
$ws.Range("B3").Value = 0.02209944751381215
$ws.Range("C3").Value = 0.02762430939226519
$ws.Range("J3").Value = 0.04972375690607735
$ws.Range("P3").Value = 0.7403314917127072
$ws.Range("S3").Value = 0.1602209944751381
$ws.Range("P4").Value = 0.7027027027027027
$ws.Range("S4").Value = 0.2972972972972973
$ws.Range("S5").Value = 1
$ws.Range("B6").Value = 0.07534246575342465
$ws.Range("D6").Value = 0.01027397260273973
$ws.Range("F6").Value = 0.1267123287671233
$ws.Range("J6").Value = 0.2226027397260274
$ws.Range("O6").Value = 0.02397260273972603
$ws.Range("Q6").Value = 0.160958904109589
$ws.Range("R6").Value = 0.07534246575342465
$ws.Range("S6").Value = 0.3047945205479452
$ws.Range("B7").Value = 0.16
$ws.Range("D7").Value = 0.035
$ws.Range("F7").Value = 0.06
$ws.Range("J7").Value = 0.13
$ws.Range("O7").Value = 0.02
$ws.Range("Q7").Value = 0.115
$ws.Range("R7").Value = 0.09
$ws.Range("S7").Value = 0.39
$ws.Range("B8").Value = 0.08597285067873303
$ws.Range("D8").Value = 0.006787330316742082
$ws.Range("F8").Value = 0.08144796380090498
$ws.Range("J8").Value = 0.1040723981900453
$ws.Range("O8").Value = 0.03619909502262444
$ws.Range("Q8").Value = 0.1719457013574661
$ws.Range("R8").Value = 0.09276018099547512
$ws.Range("S8").Value = 0.4208144796380091
$ws.Range("B9").Value = 0.1063829787234043
$ws.Range("D9").Value = 0.007092198581560284
$ws.Range("F9").Value = 0.1063829787234043
$ws.Range("J9").Value = 0.09929078014184398
$ws.Range("O9").Value = 0.03546099290780142
$ws.Range("Q9").Value = 0.1560283687943262
$ws.Range("R9").Value = 0.07801418439716312
$ws.Range("S9").Value = 0.4113475177304964
$ws.Range("B10").Value = 0.1186440677966102
$ws.Range("D10").Value = 0.02372881355932203
$ws.Range("E10").Value = 0.000847457627118644
$ws.Range("F10").Value = 0.08559322033898305
$ws.Range("J10").Value = 0.1033898305084746
$ws.Range("O10").Value = 0.03474576271186441
$ws.Range("Q10").Value = 0.1652542372881356
$ws.Range("R10").Value = 0.1016949152542373
$ws.Range("S10").Value = 0.3661016949152542
$ws.Range("G11").Value = 0.1317567567567567
$ws.Range("J11").Value = 0.09121621621621621
$ws.Range("K11").Value = 0.1824324324324324
$ws.Range("L11").Value = 0.5709459459459459
$ws.Range("S11").Value = 0.02364864864864865
$ws.Range("G12").Value = 0.7419354838709677
$ws.Range("J12").Value = 0.1612903225806452
$ws.Range("K12").Value = 0.01075268817204301
$ws.Range("L12").Value = 0.05913978494623656
$ws.Range("S12").Value = 0.02688172043010753
$ws.Range("G13").Value = 0.6458333333333334
$ws.Range("J13").Value = 0.2916666666666667
$ws.Range("S13").Value = 0.0625
$ws.Range("F15").Value = 0.03424657534246575
$ws.Range("H15").Value = 0.1438356164383562
$ws.Range("I15").Value = 0.03767123287671233
$ws.Range("J15").Value = 0.3424657534246575
$ws.Range("K15").Value = 0.08561643835616438
$ws.Range("M15").Value = 0.00684931506849315
$ws.Range("O15").Value = 0.1164383561643836
$ws.Range("S15").Value = 0.2328767123287671
$ws.Range("F16").Value = 0.005319148936170213
$ws.Range("H16").Value = 0.2287234042553191
$ws.Range("I16").Value = 0.05319148936170213
$ws.Range("J16").Value = 0.3936170212765958
$ws.Range("K16").Value = 0.101063829787234
$ws.Range("M16").Value = 0.05319148936170213
$ws.Range("N16").Value = 0.005319148936170213
$ws.Range("O16").Value = 0.05319148936170213
$ws.Range("S16").Value = 0.1063829787234043
$ws.Range("F17").Value = 0.0303030303030303
$ws.Range("H17").Value = 0.1735537190082645
$ws.Range("I17").Value = 0.07162534435261708
$ws.Range("J17").Value = 0.4297520661157025
$ws.Range("K17").Value = 0.1074380165289256
$ws.Range("M17").Value = 0.008264462809917356
$ws.Range("O17").Value = 0.07988980716253444
$ws.Range("S17").Value = 0.09917355371900827
$ws.Range("F18").Value = 0.02830188679245283
$ws.Range("H18").Value = 0.1745283018867924
$ws.Range("I18").Value = 0.05660377358490566
$ws.Range("J18").Value = 0.4386792452830189
$ws.Range("K18").Value = 0.05660377358490566
$ws.Range("M18").Value = 0.01415094339622642
$ws.Range("O18").Value = 0.1037735849056604
$ws.Range("S18").Value = 0.1273584905660377
$ws.Range("F19").Value = 0.02761982128350934
$ws.Range("H19").Value = 0.2120227457351747
$ws.Range("I19").Value = 0.06742485783915515
$ws.Range("J19").Value = 0.3403736799350122
$ws.Range("K19").Value = 0.1169780666125102
$ws.Range("M19").Value = 0.02599512591389114
$ws.Range("O19").Value = 0.08042242079610074
$ws.Range("S19").Value = 0.1291632818846466
